$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title block: "REV 0" + "3" (two runs, identical formatting) -> "REV 03"
#    (single run). Delete the trailing "3" and re-insert it right after
#    "REV 0" through a zero-length range at the same offset - the save pass
#    coalesces adjacent runs that end up with identical rPr.
# ---------------------------------------------------------------------------
$hit = $d.Content
$revFound = $hit.Find.Execute("REV 0")
if ($revFound) {
    $digitPos = $hit.End   # position right after "REV 0" = where "3" sits
    $rDigit = $d.Range($digitPos, $digitPos + 1)
    if ($rDigit.Text -eq "3") {
        $rDigit.Delete()
        $rJoin = $d.Range($digitPos, $digitPos)
        $rJoin.InsertAfter("3")
    }
}

# ---------------------------------------------------------------------------
# Locate the "HOT:" note paragraph by scanning (robust against Find leaving
# partial-range Paragraphs/Previous collections).
# ---------------------------------------------------------------------------
$paras = $d.Paragraphs
$hotIdx = -1
for ($i = 1; $i -le $paras.Count; $i++) {
    $t = $paras.Item($i).Range.Text
    if ($t -match "HOT:") {
        $hotIdx = $i
        break
    }
}

# ---------------------------------------------------------------------------
# 2) Empty paragraph right before the "HOT:" note picks up explicit
#    sz/szCs = 16 (8pt) run formatting matching its paragraph mark.
# ---------------------------------------------------------------------------
if ($hotIdx -gt 1) {
    $prevRange = $paras.Item($hotIdx - 1).Range
    if ($prevRange.Text -eq [char]13) {
        $prevRange.Font.Size = 8
        $prevRange.Font.SizeBi = 8
    }
}

# ---------------------------------------------------------------------------
# 3) Reword the "HOT:" helper note.
# ---------------------------------------------------------------------------
$oldHot = "HOT: All functions can be accessed via the menu (of course, a connected display is required). Hold UP+START+C for 1 second. Navigate the menu with A and B. Change the selected option with C and exit the menu with START."
$newHot = "HOT: All functions can be accessed through the display menu, including activating or deactivating them. Hold UP+START+C for 1 second. Navigate the menu with A and B. Change the selected option with C and exit the menu with START."
$d.Content.Find.Execute($oldHot, $true, $false, $false, $false, $false, $true, 1, $false, $newHot, 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Normal style: compatibility setting "overflow punctuation" true -> false
#    (w:overflowPunct maps to ParagraphFormat.HangingPunctuation).
# ---------------------------------------------------------------------------
$normalStyle = $d.Styles.Item("Normal")
$normalStyle.ParagraphFormat.HangingPunctuation = 0

Write-Output "edit applied"
